$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 118, pushing existing rows 118:182 down to 119:183
$ws.Rows.Item(118).Insert()

# Populate the newly inserted row 118 with the new data point
$ws.Range("A118").Value = 10
$ws.Range("B118").Value = "Vega Modelo de Temuco"
$ws.Range("C118").Value = "La Araucanía"
$ws.Range("D118").Value = 45001
$ws.Range("E118").Value = 9
$ws.Range("F118").Value = 100114002
$ws.Range("G118").Value = "Camote"
$ws.Range("H118").Value = "Sin especificar"
$ws.Range("I118").Value = "Primera"
$ws.Range("J118").Value = 50
$ws.Range("K118").Value = 26000
$ws.Range("L118").Value = 26000
$ws.Range("M118").Value = 26000
$ws.Range("N118").Value = "$/malla 20 kilos"
$ws.Range("O118").Value = "Perú"
$ws.Range("P118").Value = 1300
$ws.Range("Q118").Value = 20
$ws.Range("R118").Value = "Hortaliza"

# D118 carries the same date-number format as the rest of the date column
$ws.Range("D118").NumberFormat = $ws.Range("D119").NumberFormat
